# This script updates the weekly Fruta/Hortaliza (Breva) price records.
# The underlying source data shuffled which market-day record landed in
# which spreadsheet row (rows 2<->7, 3<->8, 6<->11, and a 5-cycle among
# rows 5/9/10/12/13), so each row is rewritten with its corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44553
$ws.Range("M2").Value = 200

# Row 3
$ws.Range("D3").Value = 44553
$ws.Range("M3").Value = 150

# Row 5
$ws.Range("D5").Value = 44189
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("R5").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S5").Value = 2143

# Row 6
$ws.Range("D6").Value = 44189
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("S6").Value = 1857

# Row 7
$ws.Range("D7").Value = 44558
$ws.Range("M7").Value = 20

# Row 8
$ws.Range("D8").Value = 44558
$ws.Range("M8").Value = 25

# Row 9
$ws.Range("D9").Value = 44187
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = "$/bandeja 7 kilos"
$ws.Range("R9").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S9").Value = 2000
$ws.Range("T9").Value = 7

# Row 10
$ws.Range("D10").Value = 44187
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 1714

# Row 11
$ws.Range("D11").Value = 44204
$ws.Range("M11").Value = 110
$ws.Range("N11").Value = 7000
$ws.Range("O11").Value = 7500
$ws.Range("P11").Value = 7318
$ws.Range("S11").Value = 1045

# Row 12
$ws.Range("D12").Value = 44550
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 24000
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 3429

# Row 13
$ws.Range("D13").Value = 44572
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = "$/bandeja 6 kilos"
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 3333
$ws.Range("T13").Value = 6

